# Auto-generated edit script applying Asura_Profits market-data refresh
# (H..N columns: currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ), LeveProfit(NQ/HQ))
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 33
$ws_ALC.Range("H33").Value = 211
$ws_ALC.Range("I33").Value = 211
$ws_ALC.Range("K33").Value = 211
$ws_ALC.Range("M33").Value = 18

# ALC row 55
$ws_ALC.Range("H55").Value = 500.5
$ws_ALC.Range("I55").Value = 500.5
$ws_ALC.Range("J55").Value = 0
$ws_ALC.Range("K55").Value = 500.5
$ws_ALC.Range("L55").Value = 0
$ws_ALC.Range("M55").ClearContents()
$ws_ALC.Range("N55").Value = -286.5

# ALC row 62
$ws_ALC.Range("H62").Value = 2043.7142
$ws_ALC.Range("I62").Value = 1533.3334
$ws_ALC.Range("J62").Value = 2426.5
$ws_ALC.Range("K62").Value = 1533.3334
$ws_ALC.Range("L62").Value = 2426.5
$ws_ALC.Range("M62").Value = -909.3334
$ws_ALC.Range("N62").Value = -3674.5

# ALC row 64
$ws_ALC.Range("H64").Value = 4105.4326
$ws_ALC.Range("I64").Value = 3144.2222
$ws_ALC.Range("J64").Value = 4414.393
$ws_ALC.Range("K64").Value = 3144.2222
$ws_ALC.Range("L64").Value = 4414.393
$ws_ALC.Range("M64").Value = -2896.2222
$ws_ALC.Range("N64").Value = -4910.393

# ALC row 65
$ws_ALC.Range("H65").Value = 2043.7142
$ws_ALC.Range("I65").Value = 1533.3334
$ws_ALC.Range("J65").Value = 2426.5
$ws_ALC.Range("K65").Value = 7666.666999999999
$ws_ALC.Range("L65").Value = 12132.5
$ws_ALC.Range("M65").Value = -4546.666999999999
$ws_ALC.Range("N65").Value = -18372.5

# ALC row 67
$ws_ALC.Range("H67").Value = 4105.4326
$ws_ALC.Range("I67").Value = 3144.2222
$ws_ALC.Range("J67").Value = 4414.393
$ws_ALC.Range("K67").Value = 3144.2222
$ws_ALC.Range("L67").Value = 4414.393
$ws_ALC.Range("M67").Value = -2286.2222
$ws_ALC.Range("N67").Value = -6130.393

# ALC row 70
$ws_ALC.Range("H70").Value = 48842.477
$ws_ALC.Range("I70").Value = 143842.86
$ws_ALC.Range("J70").Value = 1342.2858
$ws_ALC.Range("K70").Value = 431528.58
$ws_ALC.Range("L70").Value = 4026.8574
$ws_ALC.Range("M70").Value = -431258.58
$ws_ALC.Range("N70").Value = -4566.857400000001

# ALC row 73
$ws_ALC.Range("H73").Value = 48842.477
$ws_ALC.Range("I73").Value = 143842.86
$ws_ALC.Range("J73").Value = 1342.2858
$ws_ALC.Range("K73").Value = 431528.58
$ws_ALC.Range("L73").Value = 4026.8574
$ws_ALC.Range("M73").Value = -430592.58
$ws_ALC.Range("N73").Value = -5898.857400000001

# ALC row 74
$ws_ALC.Range("H74").Value = 3857.5715
$ws_ALC.Range("I74").Value = 3001.5
$ws_ALC.Range("K74").Value = 3001.5
$ws_ALC.Range("M74").Value = -2065.5

# ALC row 76
$ws_ALC.Range("H76").Value = 4337.5
$ws_ALC.Range("I76").Value = 4400
$ws_ALC.Range("K76").Value = 4400
$ws_ALC.Range("M76").Value = -4085

# ALC row 77
$ws_ALC.Range("H77").Value = 3857.5715
$ws_ALC.Range("I77").Value = 3001.5
$ws_ALC.Range("K77").Value = 15007.5
$ws_ALC.Range("M77").Value = -10327.5

# ALC row 79
$ws_ALC.Range("H79").Value = 4337.5
$ws_ALC.Range("I79").Value = 4400
$ws_ALC.Range("K79").Value = 4400
$ws_ALC.Range("M79").Value = -3308

# ALC row 129
$ws_ALC.Range("H129").Value = 968.9
$ws_ALC.Range("J129").Value = 1100.3959
$ws_ALC.Range("L129").Value = 3301.1877
$ws_ALC.Range("N129").Value = -13301.1877

# ALC row 132
$ws_ALC.Range("H132").Value = 1373.4805
$ws_ALC.Range("I132").Value = 1257.5555
$ws_ALC.Range("J132").Value = 3042.8
$ws_ALC.Range("K132").Value = 3772.6665
$ws_ALC.Range("L132").Value = 9128.400000000001
$ws_ALC.Range("M132").Value = -1242.6665
$ws_ALC.Range("N132").Value = -14188.4

# ALC row 135
$ws_ALC.Range("H135").Value = 1152.0344
$ws_ALC.Range("I135").Value = 905.95654
$ws_ALC.Range("J135").Value = 2095.3333
$ws_ALC.Range("K135").Value = 8153.60886
$ws_ALC.Range("L135").Value = 18857.9997
$ws_ALC.Range("M135").Value = -5618.60886
$ws_ALC.Range("N135").Value = -23927.9997

# ALC row 136
$ws_ALC.Range("H136").Value = 34832.5
$ws_ALC.Range("J136").Value = 34832.5
$ws_ALC.Range("L136").Value = 34832.5
$ws_ALC.Range("N136").Value = -45032.5

# ALC row 141
$ws_ALC.Range("H141").Value = 7158.1714
$ws_ALC.Range("I141").Value = 4272.64
$ws_ALC.Range("K141").Value = 12817.92
$ws_ALC.Range("M141").Value = -7637.920000000002

# ARM row 5
$ws_ARM.Range("H5").Value = 0
$ws_ARM.Range("I5").Value = 0
$ws_ARM.Range("K5").Value = 0
$ws_ARM.Range("M5").ClearContents()

# ARM row 63
$ws_ARM.Range("H63").Value = 7620.5293
$ws_ARM.Range("I63").Value = 6115
$ws_ARM.Range("K63").Value = 6115
$ws_ARM.Range("M63").Value = -5429

# ARM row 66
$ws_ARM.Range("H66").Value = 7620.5293
$ws_ARM.Range("I66").Value = 6115
$ws_ARM.Range("K66").Value = 30575
$ws_ARM.Range("M66").Value = -27143

# ARM row 88
$ws_ARM.Range("H88").Value = 2833.3333
$ws_ARM.Range("I88").Value = 1750
$ws_ARM.Range("K88").Value = 1750
$ws_ARM.Range("M88").Value = -1344

# ARM row 91
$ws_ARM.Range("H91").Value = 2833.3333
$ws_ARM.Range("I91").Value = 1750
$ws_ARM.Range("K91").Value = 1750
$ws_ARM.Range("M91").Value = -346

# BSM row 4
$ws_BSM.Range("H4").Value = 0
$ws_BSM.Range("I4").Value = 0
$ws_BSM.Range("K4").Value = 0
$ws_BSM.Range("M4").ClearContents()

# BSM row 105
$ws_BSM.Range("H105").Value = 3237.6
$ws_BSM.Range("I105").Value = 2875.111
$ws_BSM.Range("K105").Value = 2875.111
$ws_BSM.Range("M105").Value = -1128.111

# BSM row 132
$ws_BSM.Range("H132").Value = 61299.54
$ws_BSM.Range("J132").Value = 61299.54
$ws_BSM.Range("L132").Value = 61299.54
$ws_BSM.Range("N132").Value = -71419.54000000001

# CRP row 7
$ws_CRP.Range("H7").Value = 54.583332
$ws_CRP.Range("J7").Value = 81
$ws_CRP.Range("L7").Value = 81
$ws_CRP.Range("N7").Value = -307

# CRP row 62
$ws_CRP.Range("H62").Value = 67467.5
$ws_CRP.Range("I62").Value = 73962.86
$ws_CRP.Range("J62").Value = 22000
$ws_CRP.Range("K62").Value = 73962.86
$ws_CRP.Range("L62").Value = 22000
$ws_CRP.Range("M62").Value = -73338.86
$ws_CRP.Range("N62").Value = -23248

# CRP row 65
$ws_CRP.Range("H65").Value = 67467.5
$ws_CRP.Range("I65").Value = 73962.86
$ws_CRP.Range("J65").Value = 22000
$ws_CRP.Range("K65").Value = 369814.3
$ws_CRP.Range("L65").Value = 110000
$ws_CRP.Range("M65").Value = -366694.3
$ws_CRP.Range("N65").Value = -116240

# CRP row 99
$ws_CRP.Range("H99").Value = 3063.9167
$ws_CRP.Range("I99").Value = 2951.889
$ws_CRP.Range("J99").Value = 3400
$ws_CRP.Range("K99").Value = 2951.889
$ws_CRP.Range("L99").Value = 3400
$ws_CRP.Range("M99").Value = -1453.889
$ws_CRP.Range("N99").Value = -6396

# CRP row 107
$ws_CRP.Range("H107").Value = 474.81818
$ws_CRP.Range("I107").Value = 260.66666
$ws_CRP.Range("J107").Value = 653.2778
$ws_CRP.Range("K107").Value = 260.66666
$ws_CRP.Range("L107").Value = 653.2778
$ws_CRP.Range("M107").Value = 1659.33334
$ws_CRP.Range("N107").Value = -4493.2778

# CRP row 122
$ws_CRP.Range("H122").Value = 1364.2
$ws_CRP.Range("I122").Value = 1232.2667
$ws_CRP.Range("J122").Value = 1760
$ws_CRP.Range("K122").Value = 3696.800099999999
$ws_CRP.Range("L122").Value = 5280
$ws_CRP.Range("M122").Value = -1246.800099999999
$ws_CRP.Range("N122").Value = -10180

# CRP row 126
$ws_CRP.Range("H126").Value = 3063.9167
$ws_CRP.Range("I126").Value = 2951.889
$ws_CRP.Range("J126").Value = 3400
$ws_CRP.Range("K126").Value = 8855.667000000001
$ws_CRP.Range("L126").Value = 10200
$ws_CRP.Range("M126").Value = -6385.667000000001
$ws_CRP.Range("N126").Value = -15140

# CRP row 132
$ws_CRP.Range("H132").Value = 2178.2188
$ws_CRP.Range("I132").Value = 1923.7368
$ws_CRP.Range("J132").Value = 2550.1538
$ws_CRP.Range("K132").Value = 5771.2104
$ws_CRP.Range("L132").Value = 7650.4614
$ws_CRP.Range("M132").Value = -3241.2104
$ws_CRP.Range("N132").Value = -12710.4614

# CRP row 135
$ws_CRP.Range("H135").Value = 66728.57000000001
$ws_CRP.Range("J135").Value = 92688.89
$ws_CRP.Range("L135").Value = 92688.89
$ws_CRP.Range("N135").Value = -102828.89

# CUL row 23
$ws_CUL.Range("H23").Value = 344.27274
$ws_CUL.Range("I23").Value = 199.5
$ws_CUL.Range("K23").Value = 598.5
$ws_CUL.Range("M23").Value = -363.5

# CUL row 24
$ws_CUL.Range("H24").Value = 3000
$ws_CUL.Range("J24").Value = 3000
$ws_CUL.Range("L24").Value = 9000
$ws_CUL.Range("N24").Value = -9460

# CUL row 25
$ws_CUL.Range("H25").Value = 2600.1428
$ws_CUL.Range("I25").Value = 401
$ws_CUL.Range("J25").Value = 2966.6667
$ws_CUL.Range("K25").Value = 1203
$ws_CUL.Range("L25").Value = 8900.000100000001
$ws_CUL.Range("M25").Value = -1034
$ws_CUL.Range("N25").Value = -9238.000100000001

# CUL row 30
$ws_CUL.Range("H30").Value = 2600.1428
$ws_CUL.Range("I30").Value = 401
$ws_CUL.Range("J30").Value = 2966.6667
$ws_CUL.Range("K30").Value = 1203
$ws_CUL.Range("L30").Value = 8900.000100000001
$ws_CUL.Range("M30").Value = -1101
$ws_CUL.Range("N30").Value = -9104.000100000001

# CUL row 80
$ws_CUL.Range("H80").Value = 5420.2144
$ws_CUL.Range("I80").Value = 11966.667
$ws_CUL.Range("K80").Value = 35900.001
$ws_CUL.Range("M80").Value = -34964.001

# CUL row 83
$ws_CUL.Range("H83").Value = 5420.2144
$ws_CUL.Range("I83").Value = 11966.667
$ws_CUL.Range("K83").Value = 107700.003
$ws_CUL.Range("M83").Value = -103020.003

# CUL row 98
$ws_CUL.Range("H98").Value = 2254.75
$ws_CUL.Range("I98").Value = 3553.3333
$ws_CUL.Range("J98").Value = 1475.6
$ws_CUL.Range("K98").Value = 10659.9999
$ws_CUL.Range("L98").Value = 4426.799999999999
$ws_CUL.Range("M98").Value = -9161.999899999999
$ws_CUL.Range("N98").Value = -7422.799999999999

# GSM row 80
$ws_GSM.Range("H80").Value = 2913.75
$ws_GSM.Range("I80").Value = 2943.4285
$ws_GSM.Range("J80").Value = 2706
$ws_GSM.Range("K80").Value = 2943.4285
$ws_GSM.Range("L80").Value = 2706
$ws_GSM.Range("M80").Value = -1945.4285
$ws_GSM.Range("N80").Value = -4702

# GSM row 83
$ws_GSM.Range("H83").Value = 2913.75
$ws_GSM.Range("I83").Value = 2943.4285
$ws_GSM.Range("J83").Value = 2706
$ws_GSM.Range("K83").Value = 14717.1425
$ws_GSM.Range("L83").Value = 13530
$ws_GSM.Range("M83").Value = -9725.1425
$ws_GSM.Range("N83").Value = -23514

# GSM row 92
$ws_GSM.Range("H92").Value = 7472.4
$ws_GSM.Range("J92").Value = 7472.4
$ws_GSM.Range("L92").Value = 7472.4
$ws_GSM.Range("N92").Value = -11216.4

# GSM row 102
$ws_GSM.Range("H102").Value = 2481.0667
$ws_GSM.Range("I102").Value = 2441.8333
$ws_GSM.Range("K102").Value = 2441.8333
$ws_GSM.Range("M102").Value = -819.8332999999998

# GSM row 132
$ws_GSM.Range("H132").Value = 2910.8333
$ws_GSM.Range("I132").Value = 1991.5
$ws_GSM.Range("J132").Value = 4749.5
$ws_GSM.Range("K132").Value = 5974.5
$ws_GSM.Range("L132").Value = 14248.5
$ws_GSM.Range("M132").Value = -3444.5
$ws_GSM.Range("N132").Value = -19308.5

# GSM row 134
$ws_GSM.Range("H134").Value = 36732.145
$ws_GSM.Range("J134").Value = 36732.145
$ws_GSM.Range("L134").Value = 110196.435
$ws_GSM.Range("N134").Value = -115266.435

# GSM row 135
$ws_GSM.Range("H135").Value = 62753.875
$ws_GSM.Range("J135").Value = 62753.875
$ws_GSM.Range("L135").Value = 62753.875
$ws_GSM.Range("N135").Value = -72893.875

# WVR row 132
$ws_WVR.Range("H132").Value = 1902.4736
$ws_WVR.Range("I132").Value = 1262.5
$ws_WVR.Range("J132").Value = 4302.375
$ws_WVR.Range("K132").Value = 3787.5
$ws_WVR.Range("L132").Value = 12907.125
$ws_WVR.Range("M132").Value = -1257.5
$ws_WVR.Range("N132").Value = -17967.125

Write-Host "Applied Asura_Profits market data updates across ALC, ARM, BSM, CRP, CUL, GSM, WVR sheets."